$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new forecast row 54, copying the date-cell formatting from row 53 (A53)
$ws.Range("A53").Copy($ws.Range("A54"))
$ws.Range("A54").Value = 45986

$ws.Range("B54").Value = 2025
$ws.Range("C54").Value = -2.06674933094535
$ws.Range("D54").Value = 2026
$ws.Range("E54").Value = -0.3099928749133896
